{"js": "// Strike through the two \"not done yet\" todo items:\n//   \"Statistik af drift + Skrive drifttests\"\n//   \"Integrationstest + tekst til\"\n// (matches the author's \"den skal sgu lige g\u00f8res p\u00e6nere...\" cleanup commit)\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\nconst targets = [\"Statistik af drift\", \"Integrationstest + tekst til\"];\n\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const paragraph = paragraphs.items[i];\n  const text = paragraph.text;\n  if (targets.some((t) => text.indexOf(t) !== -1)) {\n    paragraph.font.strikeThrough = true;\n  }\n}\n\nawait context.sync();\n", "ps1": "# Strike through the two \"not done yet\" todo items:\n#   \"Statistik af drift + Skrive drifttests\"\n#   \"Integrationstest + tekst til\"\n# (matches the author's \"den skal sgu lige g\u00f8res p\u00e6nere...\" cleanup commit)\n\n$d = $word.ActiveDocument\n\nforeach ($p in $d.Paragraphs) {\n    $t = $p.Range.Text\n    if ($t -like \"*Statistik af drift*\" -or $t -like \"*Integrationstest + tekst til*\") {\n        $p.Range.Font.StrikeThrough = 1\n    }\n}\n"}
